$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "218.96") must be forced to
# Text format first, otherwise Excel auto-converts them to numbers instead
# of keeping the literal display string used by the source data feed.
$textForceCells = @("5:4", "8:4", "10:4", "14:4", "15:4", "18:4", "19:4", "21:4", "23:4", "24:4", "25:4", "26:4", "28:4", "29:4", "31:4", "32:4", "34:4", "36:4", "37:4", "42:4", "43:4", "44:4", "46:4", "47:4", "49:4", "50:4")
foreach ($addr in $textForceCells) {
    $parts = $addr.Split(":")
    $ws.Cells.Item([int]$parts[0], [int]$parts[1]).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "30.713.56"
$ws.Cells.Item(2, 5).Value = "  +2.76%  "
$ws.Cells.Item(3, 4).Value = "1.677.32"
$ws.Cells.Item(3, 5).Value = "  +2.89%  "
$ws.Cells.Item(4, 5).Value = "  -0.28%  "
$ws.Cells.Item(5, 4).Value = "218.96"
$ws.Cells.Item(5, 5).Value = "  +2.04%  "
$ws.Cells.Item(6, 5).Value = "  +1.98%  "
$ws.Cells.Item(7, 5).Value = "  -0.34%  "
$ws.Cells.Item(8, 4).Value = "29.10"
$ws.Cells.Item(8, 5).Value = "  +1.84%  "
$ws.Cells.Item(9, 5).Value = "  +2.16%  "
$ws.Cells.Item(10, 4).Value = "0.0643"
$ws.Cells.Item(10, 5).Value = "  +5.86%  "
$ws.Cells.Item(11, 5).Value = "  +0.07%  "
$ws.Cells.Item(12, 4).Value = "1.918.39"
$ws.Cells.Item(12, 5).Value = "  +2.92%  "
$ws.Cells.Item(13, 4).Value = "1.662.41"
$ws.Cells.Item(13, 5).Value = "  +2.04%  "
$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14, 4).Value = "10.13"
$ws.Cells.Item(14, 5).Value = "  +9.29%  "
$ws.Cells.Item(15, 2).Value = "Polygon"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(15, 4).Value = "0.605"
$ws.Cells.Item(15, 5).Value = "  +7.95%  "
$ws.Cells.Item(16, 5).Value = "  +4.69%  "
$ws.Cells.Item(17, 4).Value = "30.697.58"
$ws.Cells.Item(17, 5).Value = "  +2.62%  "
$ws.Cells.Item(18, 4).Value = "65.97"
$ws.Cells.Item(18, 5).Value = "  +3.09%  "
$ws.Cells.Item(19, 4).Value = "243.06"
$ws.Cells.Item(19, 5).Value = "  +1.06%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0719"
$ws.Cells.Item(20, 5).Value = "  +2.73%  "
$ws.Cells.Item(21, 4).Value = "0.999"
$ws.Cells.Item(21, 5).Value = "  -0.19%  "
$ws.Cells.Item(22, 5).Value = "  +2.69%  "
$ws.Cells.Item(23, 4).Value = "9.95"
$ws.Cells.Item(23, 5).Value = "  +1.79%  "
$ws.Cells.Item(24, 4).Value = "2.15"
$ws.Cells.Item(24, 5).Value = "  -0.06%  "
$ws.Cells.Item(25, 4).Value = "159.08"
$ws.Cells.Item(25, 5).Value = "  +0.76%  "
$ws.Cells.Item(26, 4).Value = "15.82"
$ws.Cells.Item(26, 5).Value = "  +2.53%  "
$ws.Cells.Item(27, 5).Value = "  +2.42%  "
$ws.Cells.Item(28, 4).Value = "6.68"
$ws.Cells.Item(28, 5).Value = "  +2.05%  "
$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 5).Value = "  -0.19%  "
$ws.Cells.Item(30, 5).Value = "  +1.06%  "
$ws.Cells.Item(31, 4).Value = "1.14"
$ws.Cells.Item(31, 5).Value = "  +3.73%  "
$ws.Cells.Item(32, 4).Value = "3.45"
$ws.Cells.Item(32, 5).Value = "  +2.67%  "
$ws.Cells.Item(33, 4).Value = "1.517.62"
$ws.Cells.Item(33, 5).Value = "  +6.60%  "
$ws.Cells.Item(34, 4).Value = "3.30"
$ws.Cells.Item(34, 5).Value = "  +4.31%  "
$ws.Cells.Item(35, 5).Value = "  +6.30%  "
$ws.Cells.Item(36, 2).Value = "TrustWalletToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(36, 4).Value = "1.02"
$ws.Cells.Item(36, 5).Value = "  -0.45%  "
$ws.Cells.Item(37, 2).Value = "Aave"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(37, 4).Value = "83.05"
$ws.Cells.Item(37, 5).Value = "  +10.39%  "
$ws.Cells.Item(38, 5).Value = "  +8.39%  "
$ws.Cells.Item(39, 5).Value = "  +4.56%  "
$ws.Cells.Item(40, 5).Value = "  -2.87%  "
$ws.Cells.Item(41, 5).Value = "  -0.09%  "
$ws.Cells.Item(42, 4).Value = "2.01"
$ws.Cells.Item(43, 4).Value = "0.836"
$ws.Cells.Item(43, 5).Value = "  +1.29%  "
$ws.Cells.Item(44, 4).Value = "0.0499"
$ws.Cells.Item(44, 5).Value = "  +0.18%  "
$ws.Cells.Item(45, 5).Value = "  +1.33%  "
$ws.Cells.Item(46, 4).Value = "0.999"
$ws.Cells.Item(46, 5).Value = "  -0.22%  "
$ws.Cells.Item(47, 4).Value = "5.54"
$ws.Cells.Item(47, 5).Value = "  +4.13%  "
$ws.Cells.Item(48, 4).Value = "1.811.77"
$ws.Cells.Item(48, 5).Value = "  +2.24%  "
$ws.Cells.Item(49, 4).Value = "49.91"
$ws.Cells.Item(49, 5).Value = "  -0.27%  "
$ws.Cells.Item(50, 4).Value = "92.63"
$ws.Cells.Item(50, 5).Value = "  +2.60%  "
$ws.Cells.Item(51, 5).Value = "  +3.23%  "
